$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K data to F:M
$ws.Range("D:E").Insert(-4161)

# Copy number formatting from the (now-shifted) F:G columns into the new D:E
# columns for each data block, so the new cells inherit the same formats
# (date format for the header rows, number format for the data rows),
# mirroring what Excel does automatically when you insert columns in the
# middle of a formatted table.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D = most recent quarter, E = prior
# quarter) with the newly reported financial figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 3895100
$ws.Range("E8").Value2 = 4030900
$ws.Range("D9").Value2 = 2727500
$ws.Range("E9").Value2 = 2714000
$ws.Range("D10").Value2 = 1167600
$ws.Range("E10").Value2 = 1316900
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 21800
$ws.Range("E14").Value2 = 20100
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 3448400
$ws.Range("E17").Value2 = 3443900
$ws.Range("D18").Value2 = 446700
$ws.Range("E18").Value2 = 587000
$ws.Range("D20").Value2 = -20400
$ws.Range("E20").Value2 = -8500
$ws.Range("D21").Value2 = 514800
$ws.Range("E21").Value2 = 663700
$ws.Range("D22").Value2 = 49000
$ws.Range("E22").Value2 = 48500
$ws.Range("D23").Value2 = 377300
$ws.Range("E23").Value2 = 530000
$ws.Range("D24").Value2 = 130400
$ws.Range("E24").Value2 = -1100
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 246900
$ws.Range("E26").Value2 = 531100
$ws.Range("D27").Value2 = 239500
$ws.Range("E27").Value2 = 526800
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 14500
$ws.Range("E29").Value2 = -11700
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 20400
$ws.Range("E32").Value2 = 8500
$ws.Range("D33").Value2 = 254000
$ws.Range("E33").Value2 = 515100
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 254000
$ws.Range("E35").Value2 = 515100
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 903400
$ws.Range("E41").Value2 = 1022500
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 2679200
$ws.Range("E43").Value2 = 2752100
$ws.Range("D44").Value2 = 1677800
$ws.Range("E44").Value2 = 1821400
$ws.Range("D45").Value2 = 471600
$ws.Range("E45").Value2 = 486600
$ws.Range("D46").Value2 = 5732000
$ws.Range("E46").Value2 = 6082600
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 1730800
$ws.Range("E48").Value2 = 1673400
$ws.Range("D49").Value2 = 9594200
$ws.Range("E49").Value2 = 9668700
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 857900
$ws.Range("E52").Value2 = 879800
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 17914900
$ws.Range("E54").Value2 = 18304500
$ws.Range("D57").Value2 = 1705300
$ws.Range("E57").Value2 = 1748700
$ws.Range("D58").Value2 = 350600
$ws.Range("E58").Value2 = 350600
$ws.Range("D59").Value2 = 2259800
$ws.Range("E59").Value2 = 2224700
$ws.Range("D60").Value2 = 4315700
$ws.Range("E60").Value2 = 4324000
$ws.Range("D61").Value2 = 3740700
$ws.Range("E61").Value2 = 3739800
$ws.Range("D62").Value2 = 2793700
$ws.Range("E62").Value2 = 2897500
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 10892200
$ws.Range("E66").Value2 = 10996400
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 9439800
$ws.Range("E72").Value2 = 9624800
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 7022700
$ws.Range("E76").Value2 = 7308100
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 254000
$ws.Range("E81").Value2 = 515100
$ws.Range("D83").Value2 = 88500
$ws.Range("E83").Value2 = 85200
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 510400
$ws.Range("E89").Value2 = 519700
$ws.Range("D91").Value2 = -114400
$ws.Range("E91").Value2 = -87800
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -108800
$ws.Range("E94").Value2 = -75700
$ws.Range("D96").Value2 = -128300
$ws.Range("E96").Value2 = -129400
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -509900
$ws.Range("E100").Value2 = -377500
$ws.Range("D101").Value2 = -10800
$ws.Range("E101").Value2 = -13500
$ws.Range("D102").Value2 = -119100
$ws.Range("E102").Value2 = 53000
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"

Write-Output "edit applied"
